$wb = $excel.ActiveWorkbook

# ---------- Sheet: Requirements Phase Defects ----------
$ws1 = $wb.Worksheets.Item("Requirements Phase Defects")
$ws1.Range("E27").Value = "2h"

# ---------- Sheet: Architect. Design Phase Defects ----------
$ws2 = $wb.Worksheets.Item("Architect. Design Phase Defects")
$ws2.Range("C10").Value = "A03"
$ws2.Range("E10").Value = "Yes, the arhitecture accounts for all the requirements"
$ws2.Range("C11").Value = "A09"
$ws2.Range("E11").Value = "relationships are not named properly"
$ws2.Range("C12").Value = "A05"
$ws2.Range("E12").Value = "No, errors are not handle"
$ws2.Range("C13").Value = "A07"
$ws2.Range("E13").Value = "No, PizzaService is PaymentService"

# ---------- Sheet: Coding Phase Defects ----------
$ws3 = $wb.Worksheets.Item("Coding Phase Defects")
$ws3.Range("D5").Value = "Dragan Alexandru"
$ws3.Range("C10").Value = "C06"
$ws3.Range("D10").Value = "MenuRepository:42"
$ws3.Range("E10").Value = "double conversion can throw exception"
$ws3.Range("C11").Value = "C07"
$ws3.Range("E11").Value = "when client(cook) click on Ready without selecting an order will throw error"
$ws3.Range("C12").Value = "C11"
$ws3.Range("D12").Value = "PizzaService:32-37"
$ws3.Range("E12").Value = "Confuzion in use of variables"
$ws3.Range("C13").Value = "C01"
$ws3.Range("D13").Value = "OrdersGuiController"
$ws3.Range("E13").Value = "* Disable payOrder until order is processed                         * Remove 0 from quantities"
$ws3.Range("E32").Value = "30min"

# Column widths on Coding Phase Defects
$ws3.Columns.Item(4).ColumnWidth = 25.7109375
$ws3.Columns.Item(5).ColumnWidth = 46.42578125

# Row heights (wrap-text rows)
$ws3.Rows.Item(11).RowHeight = 30
$ws3.Rows.Item(13).RowHeight = 30

$ws2.Rows.Item(10).RowHeight = 30

# ---------- Sheet selections / views ----------
$ws1.Range("E27").Select()
$ws2.Range("G13").Select()
$ws3.Range("E36").Select()

# Activate "Coding Phase Defects" as the active sheet/tab
$ws3.Activate()

$wb.Save()
